# Rename the "Certified deaths - X" metric labels to "# Certified deaths - X"
# and re-sort them alphabetically (they occupy rows 39-50, column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    39 = "# Certified deaths - COVID-19"
    40 = "# Certified deaths - Cancer"
    41 = "# Certified deaths - Cerebrovascular diseases"
    42 = "# Certified deaths - Chronic lower respiratory conditions"
    43 = "# Certified deaths - Dementia including Alzheimers"
    44 = "# Certified deaths - Diabetes"
    45 = "# Certified deaths - Influenza and pneumonia"
    46 = "# Certified deaths - Ischaemic heart diseases"
    47 = "# Certified deaths - Other cardiac conditions"
    48 = "# Certified deaths - Pneumonia"
    49 = "# Certified deaths - Respiratory diseases"
    50 = "# Certified deaths - Total"
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row]
}

# Match the saved selection state recorded in the workbook.
$ws.Range("C51").Select() | Out-Null
